$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Find the paragraph "Changed the app code to use inheritance"
#    It currently consists of 4 runs:
#      "Changed the " | "app" | " code to use " | "inheritance"
#    The target merges the first three runs into one run
#    ("Changed the app code to use ") while keeping "inheritance"
#    as its own, separate run.
# ------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r")
    if ($paraText -eq "Changed the app code to use inheritance") {
        $targetIndex = $i
        break
    }
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$pStart = $targetPara.Range.Start
$prefix = "Changed the app code to use "

# Force the backend to re-normalize/merge the runs inside the paragraph by
# performing a genuine text mutation (placeholder swap), then restoring the
# exact original text.
$prefixRange = $d.Range($pStart, $pStart + $prefix.Length)
$placeholder = "TEMP_PLACEHOLDER_FOR_MERGE"
$prefixRange.Text = $placeholder
$prefixRange2 = $d.Range($pStart, $pStart + $placeholder.Length)
$prefixRange2.Text = $prefix

# At this point the whole paragraph (including "inheritance") has been
# collapsed into a single run. Re-split "inheritance" back into its own
# run by toggling a character formatting property on it (same value in,
# same value out -> no visible formatting change, but forces a run break).
$targetPara = $d.Paragraphs.Item($targetIndex)
$fullText = $targetPara.Range.Text
$bodyLen = $fullText.Length - 1   # exclude the trailing paragraph mark
$suffix = "inheritance"
$suffixStart = $targetPara.Range.Start + $bodyLen - $suffix.Length
$suffixEnd = $targetPara.Range.Start + $bodyLen
$suffixRange = $d.Range($suffixStart, $suffixEnd)
$origBold = $suffixRange.Font.Bold
$suffixRange.Font.Bold = 1
$suffixRange.Font.Bold = $origBold

# ------------------------------------------------------------------
# 2) Insert the new log entries after that paragraph, and before the
#    following (originally empty) paragraph.
# ------------------------------------------------------------------
$nextIndex = $targetIndex + 1
$insertionPoint = $d.Paragraphs.Item($nextIndex).Range

$newText = "17/02/2021`r" + `
    "Sending WIFI names from ESP32 to the app`r" + `
    "18/02/2021`r" + `
    "Sending WIFI details from the app to the ESP32 method`r" + `
    "Receiving the data on the ESP32 and Connect to the WIFI`r" + `
    "X`r"

$insertionPoint.InsertBefore($newText)

# Remove the "X" placeholder character used to create a clean, run-less
# empty paragraph (avoids leaving a stray empty <w:r> behind).
$blankIndex = $targetIndex + 6
$blankPara = $d.Paragraphs.Item($blankIndex)
$xStart = $blankPara.Range.Start
$xRange = $d.Range($xStart, $xStart + 1)
$xRange.Delete()
